$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the crypto-price refresh diff.
# Column D price cells that now look like plain numbers need to stay text
# (matching the original inlineStr storage), so force NumberFormat "@" first.

$ws.Range("D2").Value = "26.846.51"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.543.50"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.07"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.41"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "1.763.27"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "1.541.43"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "26.852.49"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.32"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.53"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0683"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.06"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0458"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").Value = "1.367.87"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.960"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.74"
$ws.Range("E41").Value = "  +7.84%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.992"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").Value = "1.677.75"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.22"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Value = "0.0₇0969"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0942"
$ws.Range("E51").Value = "  -1.20%  "
